# Applies the scheduled market-data refresh to the Bahamut profit sheets.
# Updates current market-board averages (columns H-N) for the rows whose
# underlying item prices changed since the last run; values are plain
# numeric overwrites (no formulas are used in these sheets).

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 7046.5557
$ws.Range("J48").Value = 7046.5557
$ws.Range("L48").Value = 21139.6671
$ws.Range("N48").Value = -21723.6671
$ws.Range("H56").Value = 7046.5557
$ws.Range("J56").Value = 7046.5557
$ws.Range("L56").Value = 21139.6671
$ws.Range("N56").Value = -22207.6671
$ws.Range("H132").Value = 1596.2439
$ws.Range("I132").Value = 1816
$ws.Range("J132").Value = 1065.1666
$ws.Range("K132").Value = 5448
$ws.Range("L132").Value = 3195.4998
$ws.Range("M132").Value = -2918
$ws.Range("N132").Value = -8255.4998
$ws.Range("H135").Value = 1159.1765
$ws.Range("I135").Value = 658.4516
$ws.Range("J135").Value = 6333.3335
$ws.Range("K135").Value = 5926.0644
$ws.Range("L135").Value = 57000.0015
$ws.Range("M135").Value = -3391.0644
$ws.Range("N135").Value = -62070.0015
$ws.Range("H138").Value = 1072.31
$ws.Range("I138").Value = 656.42645
$ws.Range("J138").Value = 1956.0625
$ws.Range("K138").Value = 1969.27935
$ws.Range("L138").Value = 5868.1875
$ws.Range("M138").Value = 3170.72065
$ws.Range("N138").Value = -16148.1875
$ws.Range("H141").Value = 1880.3529
$ws.Range("I141").Value = 733.8946999999999
$ws.Range("J141").Value = 5231.5386
$ws.Range("K141").Value = 2201.6841
$ws.Range("L141").Value = 15694.6158
$ws.Range("M141").Value = 2978.3159
$ws.Range("N141").Value = -26054.6158

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 785.625
$ws.Range("I122").Value = 724.4167
$ws.Range("J122").Value = 969.25
$ws.Range("K122").Value = 2173.2501
$ws.Range("L122").Value = 2907.75
$ws.Range("M122").Value = 276.7498999999998
$ws.Range("N122").Value = -7807.75
$ws.Range("H123").Value = 20000
$ws.Range("J123").Value = 20000
$ws.Range("L123").Value = 20000
$ws.Range("N123").Value = -29800
$ws.Range("H132").Value = 1058.1086
$ws.Range("I132").Value = 970.4524
$ws.Range("J132").Value = 1978.5
$ws.Range("K132").Value = 2911.3572
$ws.Range("L132").Value = 5935.5
$ws.Range("M132").Value = -381.3571999999999
$ws.Range("N132").Value = -10995.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 13812
$ws.Range("I134").Value = 1137.8636
$ws.Range("J134").Value = 73561.5
$ws.Range("K134").Value = 3413.5908
$ws.Range("L134").Value = 220684.5
$ws.Range("M134").Value = -878.5907999999999
$ws.Range("N134").Value = -225754.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1157.2291
$ws.Range("I132").Value = 861.81537
$ws.Range("J132").Value = 1776.6451
$ws.Range("K132").Value = 2585.44611
$ws.Range("L132").Value = 5329.9353
$ws.Range("M132").Value = -55.44610999999986
$ws.Range("N132").Value = -10389.9353
$ws.Range("H134").Value = 992.28723
$ws.Range("I134").Value = 931.4675
$ws.Range("J134").Value = 1267.7646
$ws.Range("K134").Value = 2794.4025
$ws.Range("L134").Value = 3803.2938
$ws.Range("M134").Value = -259.4025000000001
$ws.Range("N134").Value = -8873.293799999999

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 36477.09
$ws.Range("I11").Value = 66741.336
$ws.Range("J11").Value = 160
$ws.Range("K11").Value = 200224.008
$ws.Range("L11").Value = 480
$ws.Range("M11").Value = -200084.008
$ws.Range("N11").Value = -760
$ws.Range("H12").Value = 37.46154
$ws.Range("J12").Value = 30.571428
$ws.Range("L12").Value = 91.71428400000001
$ws.Range("N12").Value = -437.714284
$ws.Range("H45").Value = 1076.2222
$ws.Range("J45").Value = 1252.2858
$ws.Range("L45").Value = 3756.8574
$ws.Range("N45").Value = -4820.857400000001
$ws.Range("H74").Value = 5386.6665
$ws.Range("J74").Value = 5812.5
$ws.Range("L74").Value = 17437.5
$ws.Range("N74").Value = -19559.5
$ws.Range("H77").Value = 5386.6665
$ws.Range("J77").Value = 5812.5
$ws.Range("L77").Value = 52312.5
$ws.Range("N77").Value = -62920.5
$ws.Range("H87").Value = 10795.5
$ws.Range("I87").Value = 1703.5
$ws.Range("J87").Value = 19887.5
$ws.Range("K87").Value = 5110.5
$ws.Range("L87").Value = 59662.5
$ws.Range("M87").Value = -3862.5
$ws.Range("N87").Value = -62158.5
$ws.Range("H90").Value = 10795.5
$ws.Range("I90").Value = 1703.5
$ws.Range("J90").Value = 19887.5
$ws.Range("K90").Value = 15331.5
$ws.Range("L90").Value = 178987.5
$ws.Range("M90").Value = -9091.5
$ws.Range("N90").Value = -191467.5
$ws.Range("H130").Value = 5000
$ws.Range("J130").Value = 12000
$ws.Range("L130").Value = 36000
$ws.Range("N130").Value = -46040
$ws.Range("H131").Value = 5568082
$ws.Range("I131").Value = 45546576
$ws.Range("J131").Value = 1456.3798
$ws.Range("K131").Value = 136639728
$ws.Range("L131").Value = 4369.1394
$ws.Range("M131").Value = -136634688
$ws.Range("N131").Value = -14449.1394

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 5136.857
$ws.Range("I53").Value = 4996.3335
$ws.Range("K53").Value = 4996.3335
$ws.Range("M53").Value = -4365.3335
$ws.Range("H70").Value = 4367.273
$ws.Range("I70").Value = 4098.3335
$ws.Range("K70").Value = 4098.3335
$ws.Range("M70").Value = -3828.3335
$ws.Range("H73").Value = 4367.273
$ws.Range("I73").Value = 4098.3335
$ws.Range("K73").Value = 4098.3335
$ws.Range("M73").Value = -3162.3335
$ws.Range("H102").Value = 1477.0667
$ws.Range("I102").Value = 1499.6666
$ws.Range("J102").Value = 1386.6666
$ws.Range("K102").Value = 1499.6666
$ws.Range("L102").Value = 1386.6666
$ws.Range("M102").Value = 122.3334
$ws.Range("N102").Value = -4630.6666
$ws.Range("H122").Value = 5645394.5
$ws.Range("I122").Value = 4989804
$ws.Range("J122").Value = 7143886.5
$ws.Range("K122").Value = 14969412
$ws.Range("L122").Value = 21431659.5
$ws.Range("M122").Value = -14966962
$ws.Range("N122").Value = -21436559.5
$ws.Range("H126").Value = 3707.8
$ws.Range("I126").Value = 3381.25
$ws.Range("J126").Value = 5014
$ws.Range("K126").Value = 10143.75
$ws.Range("L126").Value = 15042
$ws.Range("M126").Value = -7673.75
$ws.Range("N126").Value = -19982

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3586711.8
$ws.Range("I7").Value = 2536.95
$ws.Range("J7").Value = 10103393
$ws.Range("K7").Value = 2536.95
$ws.Range("L7").Value = 10103393
$ws.Range("M7").Value = -2424.95
$ws.Range("N7").Value = -10103617
$ws.Range("H40").Value = 440531.44
$ws.Range("I40").Value = 595427.9
$ws.Range("J40").Value = 1658.1666
$ws.Range("K40").Value = 595427.9
$ws.Range("L40").Value = 1658.1666
$ws.Range("M40").Value = -595291.9
$ws.Range("N40").Value = -1930.1666
$ws.Range("H46").Value = 848.5
$ws.Range("I46").Value = 848.5
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 848.5
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -660.5
$ws.Range("N46").ClearContents()
$ws.Range("H126").Value = 3586711.8
$ws.Range("I126").Value = 2536.95
$ws.Range("J126").Value = 10103393
$ws.Range("K126").Value = 7610.849999999999
$ws.Range("L126").Value = 30310179
$ws.Range("M126").Value = -5140.849999999999
$ws.Range("N126").Value = -30315119
$ws.Range("H136").Value = 1844.9215
$ws.Range("I136").Value = 1075.8914
$ws.Range("J136").Value = 8920
$ws.Range("K136").Value = 3227.6742
$ws.Range("L136").Value = 26760
$ws.Range("M136").Value = -677.6741999999999
$ws.Range("N136").Value = -31860

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 12425.286
$ws.Range("J41").Value = 12425.286
$ws.Range("L41").Value = 12425.286
$ws.Range("N41").Value = -13205.286
$ws.Range("H126").Value = 900
$ws.Range("I126").Value = 666.6667
$ws.Range("J126").Value = 1133.3334
$ws.Range("K126").Value = 2000.0001
$ws.Range("L126").Value = 3400.0002
$ws.Range("M126").Value = 469.9999
$ws.Range("N126").Value = -8340.0002
$ws.Range("H127").Value = 33714.5
$ws.Range("J127").Value = 33714.5
$ws.Range("L127").Value = 33714.5
$ws.Range("N127").Value = -43634.5
$ws.Range("H132").Value = 613.31744
$ws.Range("I132").Value = 596.62744
$ws.Range("J132").Value = 684.25
$ws.Range("K132").Value = 1789.88232
$ws.Range("L132").Value = 2052.75
$ws.Range("M132").Value = 740.1176800000001
$ws.Range("N132").Value = -7112.75
$ws.Range("H136").Value = 490.44
$ws.Range("I136").Value = 403.61905
$ws.Range("K136").Value = 1210.85715
$ws.Range("M136").Value = 1339.14285
